$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 5
$ws.Range("C9").Value = "This is Version 5 of File"

$ws.Range("C10").Select()
